# 8.3.1.2.xlsx - add a 2020 data column (K) to the indicator table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K3: header year cell (2020) -------------------------------------
# Matches the formatting of the bold "2013..2019" year header cells,
# but picks up the bigger (size 10) bold Times New Roman used for the
# title row, with the double (top+bottom) medium border already used
# by the 2018 column header (I3).
$ws.Range("I3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2020
$ws.Range("K3").Font.Size = 10

# --- K4: "Small enterprises" 2020 value -------------------------------
$ws.Range("K4").Value = 2.8218550629805335
$ws.Range("K4").Borders.Item(8).Weight = -4138
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").Font.Name = "Kyrghyz Times"
$ws.Range("K4").NumberFormat = "#,##0.0"
$ws.Range("K4").HorizontalAlignment = -4152

# --- K5: "Medium-sized enterprises" 2020 value ------------------------
# Reuses the bottom-medium border already used by the rest of row 5.
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 1.3005071159823327
$ws.Range("K5").Font.Size = 9
$ws.Range("K5").Font.Name = "Kyrghyz Times"
$ws.Range("K5").NumberFormat = "#,##0.0"
$ws.Range("K5").HorizontalAlignment = -4152

# --- Selection, matching the saved workbook cursor position -----------
$null = $ws.Range("L8").Select()
